$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.48
$ws.Range("L2").Value = 1.36
$ws.Range("M2").Value = 3.2
$ws.Range("N2").Value = 2.1
$ws.Range("O2").Value = 1.73
$ws.Range("W2").Value = 11
$ws.Range("X2").Value = 17
$ws.Range("Z2").Value = 8.5

# Row 4
$ws.Range("L4").Value = 1.29
$ws.Range("M4").Value = 3.75
$ws.Range("N4").Value = 1.97
$ws.Range("O4").Value = 1.93
$ws.Range("T4").Value = 9
$ws.Range("AH4").Value = 34
$ws.Range("AJ4").Value = 34

# Row 5
$ws.Range("G5").Value = 1.3
$ws.Range("H5").Value = 5.75
$ws.Range("I5").Value = 8
$ws.Range("N5").Value = 1.44
$ws.Range("O5").Value = 2.7
$ws.Range("P5").Value = 1.25
$ws.Range("Q5").Value = 3.75
$ws.Range("T5").Value = 9
$ws.Range("U5").Value = 7.5
$ws.Range("V5").Value = 8.5
$ws.Range("Z5").Value = 19
$ws.Range("AA5").Value = 11
$ws.Range("AB5").Value = 19

# Row 6
$ws.Range("G6").Value = 1.83
$ws.Range("H6").Value = 3.75
$ws.Range("I6").Value = 4.2
$ws.Range("N6").Value = 1.73
$ws.Range("O6").Value = 2.08
$ws.Range("U6").Value = 9.5
$ws.Range("W6").Value = 15
$ws.Range("X6").Value = 13

# Row 7
$ws.Range("G7").Value = 1.85
$ws.Range("H7").Value = 3.75
$ws.Range("I7").Value = 3.8
$ws.Range("K7").Value = 17
$ws.Range("L7").Value = 1.17
$ws.Range("M7").Value = 5
$ws.Range("N7").Value = 1.57
$ws.Range("O7").Value = 2.35
$ws.Range("P7").Value = 1.29
$ws.Range("Q7").Value = 3.5
$ws.Range("R7").Value = 1.53
$ws.Range("S7").Value = 2.38
$ws.Range("T7").Value = 10
$ws.Range("U7").Value = 11
$ws.Range("W7").Value = 17
$ws.Range("Z7").Value = 17
$ws.Range("AB7").Value = 12
$ws.Range("AC7").Value = 34
$ws.Range("AE7").Value = 15
$ws.Range("AF7").Value = 23
$ws.Range("AG7").Value = 13
$ws.Range("AI7").Value = 26
$ws.Range("AJ7").Value = 29

# Row 8
$ws.Range("N8").Value = 2.88
$ws.Range("O8").Value = 1.4

# Row 9
$ws.Range("G9").Value = 3.75
$ws.Range("H9").Value = 3
$ws.Range("J9").Value = 1.11
$ws.Range("K9").Value = 6.5
$ws.Range("T9").Value = 8.5
$ws.Range("U9").Value = 17
$ws.Range("V9").Value = 15
$ws.Range("X9").Value = 41
$ws.Range("Z9").Value = 6
$ws.Range("AB9").Value = 19
$ws.Range("AC9").Value = 67
$ws.Range("AE9").Value = 6

# Row 10
$ws.Range("H10").Value = 3.5
$ws.Range("L10").Value = 1.22
$ws.Range("M10").Value = 4
$ws.Range("N10").Value = 1.8
$ws.Range("O10").Value = 2
$ws.Range("T10").Value = 8.5
$ws.Range("AA10").Value = 7
$ws.Range("AD10").Value = 151

# Row 11
$ws.Range("G11").Value = 1.75
$ws.Range("H11").Value = 3.75
$ws.Range("T11").Value = 7.5
$ws.Range("W11").Value = 15
$ws.Range("AA11").Value = 7

# Row 14
$ws.Range("G14").Value = 2.88
$ws.Range("H14").Value = 3.5
$ws.Range("I14").Value = 2.25
$ws.Range("K14").Value = 15
$ws.Range("L14").Value = 1.18
$ws.Range("M14").Value = 4.5
$ws.Range("N14").Value = 1.65
$ws.Range("O14").Value = 2.2
$ws.Range("P14").Value = 1.3
$ws.Range("Q14").Value = 3.4
$ws.Range("R14").Value = 1.53
$ws.Range("S14").Value = 2.38
$ws.Range("T14").Value = 12
$ws.Range("U14").Value = 17
$ws.Range("W14").Value = 29
$ws.Range("Y14").Value = 26
$ws.Range("Z14").Value = 15
$ws.Range("AB14").Value = 12
$ws.Range("AD14").Value = 126
$ws.Range("AH14").Value = 23
$ws.Range("AI14").Value = 17

# Row 15
$ws.Range("J15").Value = 1.08
$ws.Range("K15").Value = 8
$ws.Range("L15").Value = 1.4
$ws.Range("M15").Value = 2.75
$ws.Range("N15").Value = 2.25
$ws.Range("O15").Value = 1.62

# Row 16
$ws.Range("G16").Value = 3.25
$ws.Range("I16").Value = 2.15
$ws.Range("J16").Value = 1.06
$ws.Range("K16").Value = 10
$ws.Range("U16").Value = 17
$ws.Range("Y16").Value = 34
$ws.Range("AG16").Value = 9
$ws.Range("AH16").Value = 19

# Row 18
$ws.Range("G18").Value = 1.85
$ws.Range("I18").Value = 4.2
$ws.Range("AC18").Value = 67

# Row 19
$ws.Range("G19").Value = 3
$ws.Range("H19").Value = 3.1
$ws.Range("J19").Value = 1.08
$ws.Range("K19").Value = 7.5
$ws.Range("L19").Value = 1.4
$ws.Range("M19").Value = 2.75
$ws.Range("N19").Value = 2.35
$ws.Range("O19").Value = 1.57
$ws.Range("P19").Value = 1.53
$ws.Range("Q19").Value = 2.38
$ws.Range("R19").Value = 2
$ws.Range("S19").Value = 1.75
$ws.Range("T19").Value = 8
$ws.Range("V19").Value = 12
$ws.Range("W19").Value = 34
$ws.Range("X19").Value = 29
$ws.Range("Y19").Value = 41
$ws.Range("Z19").Value = 7.5
$ws.Range("AB19").Value = 17
$ws.Range("AD19").Value = 401
$ws.Range("AE19").Value = 7

# Row 20
$ws.Range("N20").Value = 1.93
$ws.Range("O20").Value = 1.93
$ws.Range("R20").Value = 2.25
$ws.Range("S20").Value = 1.57
$ws.Range("T20").Value = 6
$ws.Range("AC20").Value = 81

# Row 21
$ws.Range("G21").Value = 7.5
$ws.Range("H21").Value = 4.5
$ws.Range("J21").Value = 1.05
$ws.Range("K21").Value = 11
$ws.Range("L21").Value = 1.25
$ws.Range("M21").Value = 3.75
$ws.Range("N21").Value = 1.88
$ws.Range("O21").Value = 1.98
$ws.Range("P21").Value = 1.36
$ws.Range("Q21").Value = 3
$ws.Range("R21").Value = 2.1
$ws.Range("S21").Value = 1.67
$ws.Range("Z21").Value = 11
$ws.Range("AA21").Value = 9
$ws.Range("AB21").Value = 23
$ws.Range("AC21").Value = 81
$ws.Range("AD21").Value = 501
$ws.Range("AE21").Value = 6.5
$ws.Range("AF21").Value = 6
$ws.Range("AG21").Value = 9

# Row 22
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = 3.6
$ws.Range("I22").Value = 2.25
$ws.Range("N22").Value = 2.03
$ws.Range("O22").Value = 1.83
$ws.Range("P22").Value = 1.4
$ws.Range("Q22").Value = 2.75
$ws.Range("T22").Value = 9.5
$ws.Range("W22").Value = 34
$ws.Range("AA22").Value = 7
$ws.Range("AE22").Value = 8

# Row 23
$ws.Range("K23").Value = 13

# Row 24
$ws.Range("G24").Value = 3.8
$ws.Range("I24").Value = 2.05
$ws.Range("K24").Value = 7.5
$ws.Range("AA24").Value = 6.5
$ws.Range("AB24").Value = 19
$ws.Range("AD24").Value = 501
$ws.Range("AF24").Value = 8.5

# Row 25
$ws.Range("G25").Value = 1.83
$ws.Range("H25").Value = 3.5
$ws.Range("I25").Value = 4.33
$ws.Range("J25").Value = 1.06
$ws.Range("K25").Value = 10
$ws.Range("L25").Value = 1.33
$ws.Range("M25").Value = 3.25
$ws.Range("N25").Value = 2.05
$ws.Range("O25").Value = 1.75
$ws.Range("P25").Value = 1.44
$ws.Range("Q25").Value = 2.63
$ws.Range("R25").Value = 1.83
$ws.Range("S25").Value = 1.83
$ws.Range("T25").Value = 6.5
$ws.Range("Y25").Value = 29
$ws.Range("Z25").Value = 9.5
$ws.Range("AA25").Value = 6.5
$ws.Range("AB25").Value = 17
$ws.Range("AD25").Value = 301
$ws.Range("AE25").Value = 11
$ws.Range("AF25").Value = 21
$ws.Range("AH25").Value = 41

# Row 26
$ws.Range("G26").Value = 3.1
$ws.Range("I26").Value = 2
$ws.Range("R26").Value = 1.8
$ws.Range("S26").Value = 1.95
$ws.Range("T26").Value = 10
$ws.Range("U26").Value = 17
$ws.Range("V26").Value = 12
$ws.Range("W26").Value = 34
$ws.Range("X26").Value = 26
$ws.Range("Y26").Value = 34
$ws.Range("Z26").Value = 10
$ws.Range("AE26").Value = 7.5
$ws.Range("AF26").Value = 10
$ws.Range("AG26").Value = 9
$ws.Range("AH26").Value = 19
$ws.Range("AI26").Value = 17
$ws.Range("AJ26").Value = 26

# Row 27
$ws.Range("G27").Value = 2.63
$ws.Range("H27").Value = 3.4
$ws.Range("J27").Value = 1.04
$ws.Range("K27").Value = 13
$ws.Range("L27").Value = 1.22
$ws.Range("M27").Value = 4
$ws.Range("N27").Value = 1.75
$ws.Range("O27").Value = 2.05
$ws.Range("R27").Value = 1.57
$ws.Range("S27").Value = 2.25
$ws.Range("Y27").Value = 26
$ws.Range("AA27").Value = 6.5
$ws.Range("AC27").Value = 41
$ws.Range("AD27").Value = 151
$ws.Range("AJ27").Value = 26

# Row 28
$ws.Range("L28").Value = 1.29
$ws.Range("M28").Value = 3.5
$ws.Range("N28").Value = 1.93
$ws.Range("O28").Value = 1.93
